# Apply cryptos list update (prices + 1h volume %) scraped on 2024-09-18.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D that are plain decimal numbers need to be forced back to
# Text format first, otherwise Excel auto-converts the literal price string
# (e.g. "541.58") into a numeric value, which would change the cell type from
# inline/shared string to a number - not what the source data looks like.
$ws.Range("D2").Value = "59.892.72"
$ws.Range("E2").Value = "  +1.82%  "
$ws.Range("D3").Value = "2.313.47"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.58"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.05"
$ws.Range("E6").Value = "  -1.72%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.576"
$ws.Range("E8").Value = "  -1.93%  "
$ws.Range("D9").Value = "2.312.99"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.49"
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.330"
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.36"
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").Value = "2.726.14"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("D16").Value = "59.911.19"
$ws.Range("E16").Value = "  +2.09%  "
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("D18").Value = "2.309.71"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.47"
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.07"
$ws.Range("E20").Value = "  -2.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "312.04"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.56"
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.69"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.67"
$ws.Range("E25").Value = "  +2.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.170"
$ws.Range("E26").Value = "  -1.59%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.73"
$ws.Range("E28").Value = "  -2.50%  "
$ws.Range("E29").Value = "  +3.21%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.05"
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("B31").Value = "SuiNetwork"
$ws.Range("C31").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.16"
$ws.Range("E31").Value = "  +1.65%  "
$ws.Range("D33").Value = "0.0₃0722"
$ws.Range("E33").Value = "  -1.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.81"
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("E35").Value = "  +3.28%  "
$ws.Range("E36").Value = "  -2.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.67"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.99"
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "316.30"
$ws.Range("E41").Value = "  +6.52%  "
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.51"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "136.35"
$ws.Range("E44").Value = "  -3.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.43"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0937"
$ws.Range("E46").Value = "  -2.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.563"
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.69"
$ws.Range("E48").Value = "  +2.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0490"
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("E50").Value = "  +16.90%  "
$ws.Range("E51").Value = "  -0.10%  "
